# Add annex case to createDoc util function
#
# The product sheet only had 3 sample rows (rows 2-4). This adds an "annex"
# of four more repetitions of the same 3 product templates (rows 5-16), and
# refreshes the first template's price (row 2, columns J/K) to the figure
# used by the new annex rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three product templates, keyed A/B/C, matching the existing rows 2/3/4.
# Row 2 (template A) is being repriced as part of this change: 5000/27300 -> 4500/24570.
$templateA = @{ B = 'Wieluń';  C = 'random'; D = 65489621; E = 'DB'; F = 5.5; G = 55; H = 5.46;               I = 'WA0'; J = 4500; K = 24570;               L = '["0620-img-1", "0620-img-2","0620-img-3","0620-img-4"]' }
$templateB = @{ B = 'Złoczew'; C = 'lctwo';  D = 68546513; E = 'DB'; F = 4.5; G = 45; H = 3.24;               I = 'WA0'; J = 3000; K = 9720;                L = '["0621-img-1", "0621-img-2","0621-img-3","0621-img-4"]' }
$templateC = @{ B = 'Płock';   C = 'asdasd'; D = 68546511; E = 'BK'; F = 2.5; G = 35; H = 4.5599999999999996; I = 'WB1'; J = 3500; K = 15959.999999999998; L = '["0612-img-1", "0612-img-2","0612-img-3","0612-img-4"]' }

# Apply the repriced template A values to the existing row 2.
$ws.Cells.Item(2, 10).Value = $templateA.J
$ws.Cells.Item(2, 11).Value = $templateA.K

# Rows 5-16 cycle through the same three templates, annex-style (A,B,C,A,B,C,...).
$templates = @($templateA, $templateB, $templateC)
for ($row = 5; $row -le 16; $row++) {
    $tpl = $templates[($row - 2) % 3]

    $ws.Cells.Item($row, 1).Value = $row - 1
    $ws.Cells.Item($row, 2).Value = $tpl.B
    $ws.Cells.Item($row, 3).Value = $tpl.C
    $ws.Cells.Item($row, 4).Value = $tpl.D
    $ws.Cells.Item($row, 5).Value = $tpl.E
    $ws.Cells.Item($row, 6).Value = $tpl.F
    $ws.Cells.Item($row, 7).Value = $tpl.G
    $ws.Cells.Item($row, 8).Value = $tpl.H
    $ws.Cells.Item($row, 9).Value = $tpl.I
    $ws.Cells.Item($row, 10).Value = $tpl.J
    $ws.Cells.Item($row, 11).Value = $tpl.K
    $ws.Cells.Item($row, 12).Value = $tpl.L
}

# Leave the selection the way it was left in the saved workbook.
$ws.Range("B14:L16").Select() | Out-Null
